$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "glycan"
$ws.Range("B1").Value = "binding_score"
$ws.Range("C1").Value = "monosaccharides"
$ws.Range("D1").Value = "motifs"
$ws.Range("E1").Value = "sasa"
$ws.Range("F1").Value = "flexibility"
$ws.Range("G1").Value = "has_multi_node_motifs"

$ws.Range("A2").Value = "Neu5Ac(a2-3)Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B2").Value = 1.428985662588979
$ws.Range("C2").Value = "['Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)']"
$ws.Range("D2").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E2").Value = 6.535104191608026
$ws.Range("F2").Value = 12.01774098556474
$ws.Range("G2").Value = $true

$ws.Range("A3").Value = "Neu5Ac(a2-6)Gal(b1-4)GlcNAc"
$ws.Range("B3").Value = 4.019055290693693
$ws.Range("C3").Value = "['Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-1)']"
$ws.Range("D3").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E3").Value = 8.905464473968182
$ws.Range("F3").Value = 2.229200336111317
$ws.Range("G3").Value = $true

$ws.Range("A4").Value = "Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man"
$ws.Range("B4").Value = -0.1340153345492894
$ws.Range("C4").Value = "['Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)']"
$ws.Range("D4").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E4").Value = 8.194194408653908
$ws.Range("F4").Value = 2.36214120762472
$ws.Range("G4").Value = $true

$ws.Range("A5").Value = "Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Neu5Ac(a2-3)Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B5").Value = 3.912803456406934
$ws.Range("C5").Value = "['Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)']"
$ws.Range("D5").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E5").Value = 7.098745467417885
$ws.Range("F5").Value = 13.86768279970992
$ws.Range("G5").Value = $true

$ws.Range("A6").Value = "Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B6").Value = 4.516238667748675
$ws.Range("C6").Value = "['Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)', 'Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)']"
$ws.Range("D6").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E6").Value = 14.52035013472606
$ws.Range("F6").Value = 24.36713712123527
$ws.Range("G6").Value = $true

$ws.Range("A7").Value = "Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man(a1-3)[Neu5Ac(a2-6)Gal(b1-4)GlcNAc(b1-2)Man(a1-6)]Man(b1-4)GlcNAc(b1-4)[Fuc(a1-6)]GlcNAc"
$ws.Range("B7").Value = 4.235280736096561
$ws.Range("C7").Value = "['Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)', 'Neu5Ac(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-2)']"
$ws.Range("D7").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E7").Value = 14.16856046663384
$ws.Range("F7").Value = 23.43111271398841
$ws.Range("G7").Value = $true

$ws.Range("A8").Value = "Neu5Gc(a2-6)Gal(b1-4)GlcNAc"
$ws.Range("B8").Value = 3.793897400657166
$ws.Range("C8").Value = "['Neu5Gc(a2-6)', 'Gal(b1-4)', 'GlcNAc(b1-1)']"
$ws.Range("D8").Value = "['Sia(a2-6)Gal(b1-3/4)GlcNAc']"
$ws.Range("E8").Value = 9.196050599880842
$ws.Range("F8").Value = 1.608618250959967
$ws.Range("G8").Value = $true

# Styling: header row (row 1) fully bold+bordered+centered/top; column A for data rows too
$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$colARange = $ws.Range("A2:A8")
$colARange.Font.Bold = $true
$colARange.Borders.LineStyle = 1
$colARange.HorizontalAlignment = -4108
$colARange.VerticalAlignment = -4160
